$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8254428142286088
$ws.Range("C2").Value = 0.2424580644419336
$ws.Range("E2").Value = 0.1069071643241379
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.290181421261515
$ws.Range("H2").Value = 0.4623243076471582
$ws.Range("I2").Value = 0.4265733295288108
$ws.Range("L2").Value = 0.2037260455224867
$ws.Range("N2").Value = 1.001352606265726
$ws.Range("O2").Value = 1.430463394040743
$ws.Range("B3").Value = 0.7358950841335172
$ws.Range("C3").Value = 0.2352112924915559
$ws.Range("E3").Value = 0.1075126439755465
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.2898735312658545
$ws.Range("H3").Value = 0.4656236401718417
$ws.Range("I3").Value = 0.4322234351431611
$ws.Range("L3").Value = 0.1944001279181009
$ws.Range("N3").Value = 0.9994073682988898
$ws.Range("O3").Value = 1.436479102087787
$ws.Range("B4").Value = 0.68084298382513
$ws.Range("C4").Value = 0.2307465325923914
$ws.Range("E4").Value = 0.1079436946287675
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.2899478489832674
$ws.Range("H4").Value = 0.467888091564042
$ws.Range("I4").Value = 0.4359770983620148
$ws.Range("L4").Value = 0.1887693035894671
$ws.Range("N4").Value = 0.9985864717116755
$ws.Range("O4").Value = 1.441222295312329
$ws.Range("B5").Value = 0.6583931096625975
$ws.Range("C5").Value = 0.2289234599382297
$ws.Range("E5").Value = 0.1081342678188033
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.2900442364953264
$ws.Range("H5").Value = 0.4688709000664986
$ws.Range("I5").Value = 0.4375782141929054
$ws.Range("L5").Value = 0.1864988012629709
$ws.Range("N5").Value = 0.998346195654932
$ws.Range("O5").Value = 1.443418865170415
$ws.Range("B6").Value = 0.6546644290600909
$ws.Range("C6").Value = 0.2286205251194247
$ws.Range("E6").Value = 0.1081668135437575
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.2900642303939804
$ws.Range("H6").Value = 0.4690377205828042
$ws.Range("I6").Value = 0.4378483927817101
$ws.Range("L6").Value = 0.1861232458525706
$ws.Range("N6").Value = 0.9983119999918131
$ws.Range("O6").Value = 1.44379952275753
$ws.Range("B7").Value = 0.6805402781877774
$ws.Range("C7").Value = 0.2307219605248605
$ws.Range("E7").Value = 0.1079462043565158
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.2899488814155049
$ws.Range("H7").Value = 0.467901102997935
$ws.Range("I7").Value = 0.435998402288039
$ws.Range("L7").Value = 0.1887385850260159
$ws.Range("N7").Value = 0.9985828492263664
$ws.Range("O7").Value = 1.441250851698442
$ws.Range("B8").Value = 0.7945821969410076
$ws.Range("C8").Value = 0.2399626396640997
$ws.Range("E8").Value = 0.107103634782673
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.2900205251478596
$ws.Range("H8").Value = 0.4634123956227683
$ws.Range("I8").Value = 0.4284624037269573
$ws.Range("L8").Value = 0.2004907609972122
$ws.Range("N8").Value = 1.00060456459012
$ws.Range("O8").Value = 1.432319639797782
$ws.Range("B9").Value = 1.017600340701335
$ws.Range("C9").Value = 0.25795587458461
$ws.Range("E9").Value = 0.1059214076163464
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.2922567402804717
$ws.Range("H9").Value = 0.4565030759508701
$ws.Range("I9").Value = 0.4159444246601289
$ws.Range("L9").Value = 0.2242885843422897
$ws.Range("N9").Value = 1.007519421790278
$ws.Range("O9").Value = 1.423145016604025
$ws.Range("B10").Value = 1.181001092253439
$ws.Range("C10").Value = 0.27108947553549
$ws.Range("E10").Value = 0.1053390157484628
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.2951864307365639
$ws.Range("H10").Value = 0.4525803608331813
$ws.Range("I10").Value = 0.4081291899397925
$ws.Range("L10").Value = 0.2422270677968754
$ws.Range("N10").Value = 1.014382536456068
$ws.Range("O10").Value = 1.421507109936982
$ws.Range("B11").Value = 1.255223710777045
$ws.Range("C11").Value = 0.2770439567402434
$ws.Range("E11").Value = 0.1051361427356952
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.2968006193682129
$ws.Range("H11").Value = 0.4510462231733641
$ws.Range("I11").Value = 0.4048747020259214
$ws.Range("L11").Value = 0.2504856170171905
$ws.Range("N11").Value = 1.017888700603081
$ws.Range("O11").Value = 1.42187412167678
$ws.Range("B12").Value = 1.283312591610184
$ws.Range("C12").Value = 0.2792957138318855
$ws.Range("E12").Value = 0.1050682367252485
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2974524871383295
$ws.Range("H12").Value = 0.450501273157073
$ws.Range("I12").Value = 0.403685626053182
$ws.Range("L12").Value = 0.25362693204859
$ws.Range("N12").Value = 1.019271339359207
$ws.Range("O12").Value = 1.422173310440286
$ws.Range("B13").Value = 1.277263962317932
$ws.Range("C13").Value = 0.2788108974822308
$ws.Range("E13").Value = 0.1050824649929041
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.2973102874758524
$ws.Range("H13").Value = 0.4506170370538314
$ws.Range("I13").Value = 0.4039397864298628
$ws.Range("L13").Value = 0.2529497741368942
$ws.Range("N13").Value = 1.01897112508577
$ws.Range("O13").Value = 1.422101744849584
$ws.Range("B14").Value = 1.257534963850674
$ws.Range("C14").Value = 0.2772292728867285
$ws.Range("E14").Value = 0.1051303773834498
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.2968534342707017
$ws.Range("H14").Value = 0.4510006684325845
$ws.Range("I14").Value = 0.4047760071740925
$ws.Range("L14").Value = 0.2507437755174777
$ws.Range("N14").Value = 1.01800135195397
$ws.Range("O14").Value = 1.421895523632372
$ws.Range("B15").Value = 1.245448032121146
$ws.Range("C15").Value = 0.2762600759422469
$ws.Range("E15").Value = 0.1051608862438371
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.2965788910740486
$ws.Range("H15").Value = 0.4512403412342536
$ws.Range("I15").Value = 0.4052938616916535
$ws.Range("L15").Value = 0.249394354163087
$ws.Range("N15").Value = 1.017414482577976
$ws.Range("O15").Value = 1.421790079723962
$ws.Range("B16").Value = 1.176148108719133
$ws.Range("C16").Value = 0.2706999161765964
$ws.Range("E16").Value = 0.1053535222250908
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.2950866150484046
$ws.Range("H16").Value = 0.4526856573652509
$ws.Range("I16").Value = 0.4083479367513547
$ws.Range("L16").Value = 0.2416893166557657
$ws.Range("N16").Value = 1.014161103478102
$ws.Range("O16").Value = 1.421505528119127
$ws.Range("B17").Value = 1.133605520182584
$ws.Range("C17").Value = 0.2672836642560412
$ws.Range("E17").Value = 0.1054875885853015
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.2942433375851365
$ws.Range("H17").Value = 0.4536364250167537
$ws.Range("I17").Value = 0.410298592161503
$ws.Range("L17").Value = 0.2369875847471121
$ws.Range("N17").Value = 1.01226342109932
$ws.Range("O17").Value = 1.421616009021676
$ws.Range("B18").Value = 1.109125983746992
$ws.Range("C18").Value = 0.2653168497737681
$ws.Range("E18").Value = 0.1055705417880208
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.2937847907780053
$ws.Range("H18").Value = 0.4542068433249753
$ws.Range("I18").Value = 0.4114488558214795
$ws.Range("L18").Value = 0.2342925287576776
$ws.Range("N18").Value = 1.011208095386323
$ws.Range("O18").Value = 1.421784214056316
$ws.Range("B19").Value = 1.100835945807376
$ws.Range("C19").Value = 0.2646506036890059
$ws.Range("E19").Value = 0.1055996318456547
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.2936340791976093
$ws.Range("H19").Value = 0.4544040239968865
$ws.Range("I19").Value = 0.4118431729127607
$ws.Range("L19").Value = 0.2333816232334271
$ws.Range("N19").Value = 1.010857001926041
$ws.Range("O19").Value = 1.42185913183252
$ws.Range("B20").Value = 1.138135317218655
$ws.Range("C20").Value = 0.2676475257266873
$ws.Range("E20").Value = 0.105472712437475
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.294330363915023
$ws.Range("H20").Value = 0.4535327756209
$ws.Range("I20").Value = 0.4100880121639179
$ws.Range("L20").Value = 0.2374871352426737
$ws.Range("N20").Value = 1.012461691270119
$ws.Range("O20").Value = 1.421593414177863
$ws.Range("B21").Value = 1.263330342185725
$ws.Range("C21").Value = 0.2776939193635144
$ws.Range("E21").Value = 0.1051160623846989
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.2969865199209352
$ws.Range("H21").Value = 0.4508870097292288
$ws.Range("I21").Value = 0.4045292122684785
$ws.Range("L21").Value = 0.2513913524906144
$ws.Range("N21").Value = 1.018284709421351
$ws.Range("O21").Value = 1.421951745471517
$ws.Range("B22").Value = 1.345048976484748
$ws.Range("C22").Value = 0.284241802528129
$ws.Range("E22").Value = 0.1049349464037164
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2989592371948504
$ws.Range("H22").Value = 0.4493676520067709
$ws.Range("I22").Value = 0.4011488039331006
$ws.Range("L22").Value = 0.2605599901648503
$ws.Range("N22").Value = 1.022410428049298
$ws.Range("O22").Value = 1.423119890642511
$ws.Range("B23").Value = 1.301444308800058
$ws.Range("C23").Value = 0.2807487854863666
$ws.Range("E23").Value = 0.1050268578519677
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.2978846521488805
$ws.Range("H23").Value = 0.4501593654523077
$ws.Range("I23").Value = 0.4029298508464656
$ws.Range("L23").Value = 0.2556591132747741
$ws.Range("N23").Value = 1.020179269524519
$ws.Range("O23").Value = 1.422410877700941
$ws.Range("B24").Value = 1.136087461977752
$ws.Range("C24").Value = 0.2674830325556172
$ws.Range("E24").Value = 0.1054794196384918
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.2942909375101408
$ws.Range("H24").Value = 0.4535795613478797
$ws.Range("I24").Value = 0.4101831256595254
$ws.Range("L24").Value = 0.2372612634341635
$ws.Range("N24").Value = 1.012371942171129
$ws.Range("O24").Value = 1.421603303234178
$ws.Range("B25").Value = 0.9573418291083726
$ws.Range("C25").Value = 0.2531027473726368
$ws.Range("E25").Value = 0.1061909466462136
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.2914264741208825
$ws.Range("H25").Value = 0.4581696178555106
$ws.Range("I25").Value = 0.419088505419321
$ws.Range("L25").Value = 0.2177705767815468
$ws.Range("N25").Value = 1.005334573183177
$ws.Range("O25").Value = 1.424732251527175
